# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx price/volume/coin-ranking update described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.731.55"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "3.400.12"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'408.27"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "'127.38"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("D7").Value = "'0.614"
$ws.Range("E7").Value = "  -2.56%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("E10").Value = "  -9.69%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "3.942.95"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.140"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "'9.01"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("E15").Value = "  -8.38%  "
$ws.Range("D16").Value = "'20.31"
$ws.Range("E16").Value = "  -4.38%  "
$ws.Range("D17").Value = "3.415.38"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Value = "'1.06"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'12.11"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").Value = "61.703.98"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").Value = "'479.57"
$ws.Range("E21").Value = "  +22.52%  "
$ws.Range("D22").Value = "'89.07"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").Value = "'33.21"
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").Value = "'9.08"
$ws.Range("E27").Value = "  +3.51%  "
$ws.Range("D28").Value = "'4.80"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'7.76"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("D31").Value = "'11.72"
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("E32").Value = "  -3.54%  "
$ws.Range("E33").Value = "  -6.48%  "
$ws.Range("D34").Value = "'40.75"
$ws.Range("E34").Value = "  -7.40%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "'56.56"
$ws.Range("E36").Value = "  +5.35%  "
$ws.Range("D37").Value = "'0.0480"
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").Value = "'149.43"
$ws.Range("E39").Value = "  +5.58%  "
$ws.Range("D40").Value = "'3.33"
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").Value = "'0.317"
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").Value = "'2.90"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").Value = "'2.05"
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'4.13"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'2.51"
$ws.Range("E46").Value = "  +2.69%  "
$ws.Range("D47").Value = "'2.31"
$ws.Range("E47").Value = "  +16.83%  "
$ws.Range("D48").Value = "'16.15"
$ws.Range("E48").Value = "  -3.82%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.144"
$ws.Range("E49").Value = "  +10.04%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'21.87"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("D51").Value = "'113.27"
$ws.Range("E51").Value = "  +16.79%  "
